$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.2280973941087723
$ws.Range("B1").Value = 0.2552430927753448
$ws.Range("C1").Value = 0.3115582168102264
$ws.Range("D1").Value = 0.5969285368919373
$ws.Range("E1").Value = 4.431256294250488
